# Add a new slide ("Link Configurations Part 2") at the end of the deck,
# reproducing plan/challenges.pptx's slide41.xml.
#
# EMU -> point conversion factor used throughout (1 pt = 12700 EMU), since
# the Shapes.Add* family of COM calls takes coordinates in points.
$EMU = 12700

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. New slide, appended after the last slide, "Title and Content" layout
#    (layout index 2) -- matches every other titled slide in this deck.
# ---------------------------------------------------------------------
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# Title placeholder text.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Link Configurations Part 2"

# The "Title and Content" layout auto-creates a Content Placeholder as the
# 2nd shape; the target slide only has the title placeholder, so drop it.
if ($s.Shapes.Count -ge 2) {
    $s.Shapes.Item(2).Delete()
}

# ---------------------------------------------------------------------
# 2. Straight connectors (thin horizontal lines).
# ---------------------------------------------------------------------
function Add-Connector($x, $y, $cx, $colorHex) {
    $x1 = $x / $EMU
    $y1 = $y / $EMU
    $x2 = ($x + $cx) / $EMU
    $y2 = $y1
    $ln = $s.Shapes.AddLine($x1, $y1, $x2, $y2)
    # AddLine can leave a sub-EMU rounding error on the (identical)
    # begin/end Y coordinate; pin the height back to exactly 0.
    $ln.Height = 0
    if ($colorHex) {
        # .RGB is an OLE_COLOR (0xBBGGRR); swap our RRGGBB hex string
        # around before handing it over so the saved <a:srgbClr val=.../>
        # ends up as the intended RRGGBB value.
        $r = $colorHex.Substring(0, 2)
        $g = $colorHex.Substring(2, 2)
        $b = $colorHex.Substring(4, 2)
        $bgr = "$b$g$r"
        $ln.Line.ForeColor.RGB = [Convert]::ToInt32($bgr, 16)
    }
    return $ln
}

Add-Connector 1849869 3948354 2081630 $null       | Out-Null
Add-Connector 1849869 4189244 2081630 $null       | Out-Null
Add-Connector 6237163 3026228 2081630 "00B050"    | Out-Null
Add-Connector 6237163 3267118 2081630 "00B050"    | Out-Null
Add-Connector 6237163 3532589 2081630 "00B050"    | Out-Null
Add-Connector 6237163 3773479 2081630 "00B050"    | Out-Null
Add-Connector 6237163 4055806 2081630 "00B050"    | Out-Null
Add-Connector 6237163 4296696 2081630 "00B050"    | Out-Null
Add-Connector 1849869 3696226 2081630 $null       | Out-Null

# ---------------------------------------------------------------------
# 3. Explanatory text box.
# ---------------------------------------------------------------------
$tbX = 4395018 / $EMU
$tbY = 5098728 / $EMU
$tbW = 5477973 / $EMU
$tbH = 1477328 / $EMU

$tb = $s.Shapes.AddTextbox(1, $tbX, $tbY, $tbW, $tbH)
$tr = $tb.TextFrame.TextRange
$tr.Text = "Link configurations for incoming lanes to outgoing lanes can be simplified if we order the outgoing lanes in such a way that leftmost outgoing lane is at the top and rightmost at the bottom. In addition, we can split, but "
$tr.InsertAfter("cannot merge.") | Out-Null

$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0

# AutoSize recomputes the box height from its content; restore the exact
# recorded size/position afterwards.
$tb.Left = $tbX
$tb.Top = $tbY
$tb.Width = $tbW
$tb.Height = $tbH

Write-Host "Slide count now: $($p.Slides.Count)"
